# Scheduled-runner market data refresh: update currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ (columns H-N) for the leves whose market
# snapshot changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 565.1053000000001
$ws.Range("I33").Value = 265.75
$ws.Range("J33").Value = 2161.6667
$ws.Range("K33").Value = 265.75
$ws.Range("L33").Value = 2161.6667
$ws.Range("M33").Value = -36.75
$ws.Range("N33").Value = -2619.6667

$ws.Range("H43").Value = 2132.6667
$ws.Range("I43").Value = 1959.2
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 1959.2
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -1890.2
$ws.Range("N43").Value = -3138

$ws.Range("H62").Value = 5198.154
$ws.Range("I62").Value = 4572.375
$ws.Range("J62").Value = 6199.4
$ws.Range("K62").Value = 4572.375
$ws.Range("L62").Value = 6199.4
$ws.Range("M62").Value = -3948.375
$ws.Range("N62").Value = -7447.4

$ws.Range("H65").Value = 5198.154
$ws.Range("I65").Value = 4572.375
$ws.Range("J65").Value = 6199.4
$ws.Range("K65").Value = 22861.875
$ws.Range("L65").Value = 30997
$ws.Range("M65").Value = -19741.875
$ws.Range("N65").Value = -37237

$ws.Range("H98").Value = 1406.8125
$ws.Range("I98").Value = 700.73334
$ws.Range("J98").Value = 11998
$ws.Range("K98").Value = 700.73334
$ws.Range("L98").Value = 11998
$ws.Range("M98").Value = 797.26666
$ws.Range("N98").Value = -14994

$ws.Range("H112").Value = 3013.4915
$ws.Range("I112").Value = 1999
$ws.Range("J112").Value = 3030.9827
$ws.Range("K112").Value = 5997
$ws.Range("L112").Value = 9092.9481
$ws.Range("M112").Value = -4889
$ws.Range("N112").Value = -11308.9481

$ws.Range("H113").Value = 2871.7273
$ws.Range("I113").Value = 1973.75
$ws.Range("J113").Value = 5266.3335
$ws.Range("K113").Value = 1973.75
$ws.Range("L113").Value = 5266.3335
$ws.Range("M113").Value = 1280.25
$ws.Range("N113").Value = -11774.3335

$ws.Range("H122").Value = 1406.8125
$ws.Range("I122").Value = 700.73334
$ws.Range("J122").Value = 11998
$ws.Range("K122").Value = 2102.20002
$ws.Range("L122").Value = 35994
$ws.Range("M122").Value = 347.7999799999998
$ws.Range("N122").Value = -40894

$ws.Range("H141").Value = 3042.6667
$ws.Range("I141").Value = 2898.1667
$ws.Range("J141").Value = 3331.6667
$ws.Range("K141").Value = 8694.500100000001
$ws.Range("L141").Value = 9995.000100000001
$ws.Range("M141").Value = -3514.500100000001
$ws.Range("N141").Value = -20355.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H32").Value = 4892.22
$ws.Range("I32").Value = 1839.7906
$ws.Range("J32").Value = 23642.857
$ws.Range("K32").Value = 1839.7906
$ws.Range("L32").Value = 23642.857
$ws.Range("M32").Value = -1552.7906
$ws.Range("N32").Value = -24216.857

$ws.Range("H61").Value = 5045
$ws.Range("I61").Value = 3879.8948
$ws.Range("J61").Value = 8207.429
$ws.Range("K61").Value = 3879.8948
$ws.Range("L61").Value = 8207.429
$ws.Range("M61").Value = -3667.8948
$ws.Range("N61").Value = -8631.429

$ws.Range("H122").Value = 4608.816
$ws.Range("I122").Value = 3997.9312
$ws.Range("J122").Value = 6577.222
$ws.Range("K122").Value = 11993.7936
$ws.Range("L122").Value = 19731.666
$ws.Range("M122").Value = -9543.793600000001
$ws.Range("N122").Value = -24631.666

$ws.Range("H132").Value = 1476.2391
$ws.Range("I132").Value = 1300.7675
$ws.Range("J132").Value = 3991.3333
$ws.Range("K132").Value = 3902.3025
$ws.Range("L132").Value = 11973.9999
$ws.Range("M132").Value = -1372.3025
$ws.Range("N132").Value = -17033.9999

$ws.Range("H136").Value = 5045
$ws.Range("I136").Value = 3879.8948
$ws.Range("J136").Value = 8207.429
$ws.Range("K136").Value = 11639.6844
$ws.Range("L136").Value = 24622.287
$ws.Range("M136").Value = -9089.6844
$ws.Range("N136").Value = -29722.287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2330.158
$ws.Range("I86").Value = 2181.25
$ws.Range("J86").Value = 2585.4285
$ws.Range("K86").Value = 2181.25
$ws.Range("L86").Value = 2585.4285
$ws.Range("M86").Value = -1058.25
$ws.Range("N86").Value = -4831.4285

$ws.Range("H89").Value = 2330.158
$ws.Range("I89").Value = 2181.25
$ws.Range("J89").Value = 2585.4285
$ws.Range("K89").Value = 10906.25
$ws.Range("L89").Value = 12927.1425
$ws.Range("M89").Value = -5290.25
$ws.Range("N89").Value = -24159.1425

$ws.Range("H92").Value = 29800.334
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 29800.334
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 29800.334
$ws.Range("N92").Value = -34792.334

$ws.Range("H105").Value = 3757.6924
$ws.Range("I105").Value = 2910.5557
$ws.Range("J105").Value = 5663.75
$ws.Range("K105").Value = 2910.5557
$ws.Range("L105").Value = 5663.75
$ws.Range("M105").Value = -1163.5557
$ws.Range("N105").Value = -9157.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6688.564
$ws.Range("I31").Value = 7484
$ws.Range("J31").Value = 6290.846
$ws.Range("K31").Value = 7484
$ws.Range("L31").Value = 6290.846
$ws.Range("M31").Value = -7189
$ws.Range("N31").Value = -6880.846

$ws.Range("H34").Value = 6688.564
$ws.Range("I34").Value = 7484
$ws.Range("J34").Value = 6290.846
$ws.Range("K34").Value = 7484
$ws.Range("L34").Value = 6290.846
$ws.Range("M34").Value = -7282
$ws.Range("N34").Value = -6694.846

$ws.Range("H99").Value = 7426.294
$ws.Range("I99").Value = 4420.6665
$ws.Range("J99").Value = 10807.625
$ws.Range("K99").Value = 4420.6665
$ws.Range("L99").Value = 10807.625
$ws.Range("M99").Value = -2922.6665
$ws.Range("N99").Value = -13803.625

$ws.Range("H107").Value = 27779464
$ws.Range("I107").Value = 38463150
$ws.Range("J107").Value = 1874
$ws.Range("K107").Value = 38463150
$ws.Range("L107").Value = 1874
$ws.Range("M107").Value = -38461230
$ws.Range("N107").Value = -5714

$ws.Range("H126").Value = 7426.294
$ws.Range("I126").Value = 4420.6665
$ws.Range("J126").Value = 10807.625
$ws.Range("K126").Value = 13261.9995
$ws.Range("L126").Value = 32422.875
$ws.Range("M126").Value = -10791.9995
$ws.Range("N126").Value = -37362.875

$ws.Range("H132").Value = 2434.5
$ws.Range("I132").Value = 2413.6667
$ws.Range("J132").Value = 2497
$ws.Range("K132").Value = 7241.000100000001
$ws.Range("L132").Value = 7491
$ws.Range("M132").Value = -4711.000100000001
$ws.Range("N132").Value = -12551

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 379.9091
$ws.Range("I33").Value = 318.42856
$ws.Range("J33").Value = 487.5
$ws.Range("K33").Value = 1910.57136
$ws.Range("L33").Value = 2925
$ws.Range("M33").Value = -1627.57136
$ws.Range("N33").Value = -3491

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 242.14285
$ws.Range("I2").Value = 287
$ws.Range("J2").Value = 130
$ws.Range("K2").Value = 287
$ws.Range("L2").Value = 130
$ws.Range("M2").Value = -174
$ws.Range("N2").Value = -356

$ws.Range("H97").Value = 385.86957
$ws.Range("I97").Value = 357.95456
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 357.95456
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 138.04544
$ws.Range("N97").Value = -1992

$ws.Range("H126").Value = 5940.478
$ws.Range("I126").Value = 5537.0713
$ws.Range("J126").Value = 6568
$ws.Range("K126").Value = 16611.2139
$ws.Range("L126").Value = 19704
$ws.Range("M126").Value = -14141.2139
$ws.Range("N126").Value = -24644

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4877
$ws.Range("I40").Value = 4722.517
$ws.Range("J40").Value = 6370.3335
$ws.Range("K40").Value = 4722.517
$ws.Range("L40").Value = 6370.3335
$ws.Range("M40").Value = -4586.517
$ws.Range("N40").Value = -6642.3335

$ws.Range("H68").Value = 4321.1665
$ws.Range("I68").Value = 3983.75
$ws.Range("J68").Value = 4996
$ws.Range("K68").Value = 3983.75
$ws.Range("L68").Value = 4996
$ws.Range("M68").Value = -3234.75
$ws.Range("N68").Value = -6494

$ws.Range("H71").Value = 4321.1665
$ws.Range("I71").Value = 3983.75
$ws.Range("J71").Value = 4996
$ws.Range("K71").Value = 19918.75
$ws.Range("L71").Value = 24980
$ws.Range("M71").Value = -16174.75
$ws.Range("N71").Value = -32468

$ws.Range("H82").Value = 1764.5
$ws.Range("I82").Value = 1689.1666
$ws.Range("J82").Value = 1990.5
$ws.Range("K82").Value = 1689.1666
$ws.Range("L82").Value = 1990.5
$ws.Range("M82").Value = -1328.1666
$ws.Range("N82").Value = -2712.5

$ws.Range("H85").Value = 1764.5
$ws.Range("I85").Value = 1689.1666
$ws.Range("J85").Value = 1990.5
$ws.Range("K85").Value = 1689.1666
$ws.Range("L85").Value = 1990.5
$ws.Range("M85").Value = -441.1666
$ws.Range("N85").Value = -4486.5

$ws.Range("H132").Value = 5549.7856
$ws.Range("I132").Value = 6321.8696
$ws.Range("J132").Value = 1998.2
$ws.Range("K132").Value = 18965.6088
$ws.Range("L132").Value = 5994.6
$ws.Range("M132").Value = -16435.6088
$ws.Range("N132").Value = -11054.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3195.0625
$ws.Range("I122").Value = 2294
$ws.Range("J122").Value = 5177.4
$ws.Range("K122").Value = 6882
$ws.Range("L122").Value = 15532.2
$ws.Range("M122").Value = -4432
$ws.Range("N122").Value = -20432.2

$ws.Range("H126").Value = 3863.0715
$ws.Range("I126").Value = 3138.5
$ws.Range("J126").Value = 5674.5
$ws.Range("K126").Value = 9415.5
$ws.Range("L126").Value = 17023.5
$ws.Range("M126").Value = -6945.5
$ws.Range("N126").Value = -21963.5

$ws.Range("H132").Value = 5356.2163
$ws.Range("I132").Value = 3068.2812
$ws.Range("J132").Value = 19999
$ws.Range("K132").Value = 9204.8436
$ws.Range("L132").Value = 59997
$ws.Range("M132").Value = -6674.8436
$ws.Range("N132").Value = -65057
